$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 29 - Home module/component
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A29:D29").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(29, 1).Value = 44987
$ws.Cells.Item(29, 2).Value = 26
$ws.Cells.Item(29, 3).Value = "Create Home Module and Home Component"
$ws.Cells.Item(29, 4).Value = "ng g m ./pages/home`nng g c ./pages/home --skip-tests"
$ws.Rows.Item(29).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 30 - Performance module/component
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A30:D30").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(30, 1).Value = 44987
$ws.Cells.Item(30, 2).Value = 27
$ws.Cells.Item(30, 3).Value = "Create Performance Module and Performance Component"
$ws.Cells.Item(30, 4).Value = "ng g m ./pages/performance`nng g c ./pages/performance --skip-tests"
$ws.Rows.Item(30).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 31 - SEP module/component
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A31:D31").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(31, 1).Value = 44987
$ws.Cells.Item(31, 2).Value = 28
$ws.Cells.Item(31, 3).Value = "Create SEP Module and SEP Component"
$ws.Cells.Item(31, 4).Value = "ng g m ./pages/sep`nng g c ./pages/sep --skip-tests"
$ws.Rows.Item(31).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 32 - E-TS1 module/component
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A32:D32").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(32, 1).Value = 44987
$ws.Cells.Item(32, 2).Value = 29
$ws.Cells.Item(32, 3).Value = "Create E-TS1 Module and E-TS1 Component"
$ws.Cells.Item(32, 4).Value = "ng g m ./pages/ets1 --routing`nng g c ./pages/ets1 --skip-tests`nng g c ./pages/ets1/iworkspace --skip-tests`nng g c ./pages/ets1/tworkspace --skip-tests"
$ws.Rows.Item(32).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 33 - switch to feature/sep-menu branch
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A33:D33").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(33, 1).Value = 44994
$ws.Cells.Item(33, 2).Value = 30
$ws.Cells.Item(33, 3).Value = "Switch to the new branch by using 'git switch -c ""feature/sep-menu""'"
$ws.Cells.Item(33, 4).Value = 'git switch -c "feature/sep-menu"'

# ---------------------------------------------------------------------------
# Row 34 - SEP card component (Detail cell reuses the "C" style (s=1)
# instead of the usual Detail style, and gets its own switch/commit
# highlight conditional format, same colors as the existing C/E rules).
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A34:D34").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(34, 1).Value = 44994
$ws.Cells.Item(34, 2).Value = 31
$ws.Cells.Item(34, 3).Value = 'Create SEP Card Componenet by using "ng g c ./pages/one-sep-card --skip-tests"'

$ws.Range("C28").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(34, 4).Value = "ng g c ./pages/one-sep-card --skip-tests"

$d34 = $ws.Range("D34")
$ruleSwitch = $d34.FormatConditions.Add(9, $null, $null, $null, "switch", 0)
$ruleSwitch.Font.Color = 24832
$ruleSwitch.Interior.Color = 13561798
$ruleSwitch.Priority = 1
$ruleCommit = $d34.FormatConditions.Add(9, $null, $null, $null, "commit", 0)
$ruleCommit.Font.Color = 22428
$ruleCommit.Interior.Color = 10284031
$ruleCommit.Priority = 2

# Existing C / E conditional formats get pushed down in priority now that
# the new D34 rules take priority 1-2.
$cConds = $ws.Range("C1:C1048576").FormatConditions
$cConds.Item(1).Priority = 5
$cConds.Item(2).Priority = 6
$eConds = $ws.Range("E2:E1048576").FormatConditions
$eConds.Item(1).Priority = 3
$eConds.Item(2).Priority = 4

# ---------------------------------------------------------------------------
# Row 35 - SEP card interface
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A35:D35").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(35, 1).Value = 44994
$ws.Cells.Item(35, 2).Value = 32
$ws.Cells.Item(35, 3).Value = 'Create SEP Card Interface by using "ng g i ./@core/shared/interfaces/one-sep-card"'
$ws.Cells.Item(35, 4).Value = "ng g i ./@core/shared/interfaces/one-sep-card"
$ws.Rows.Item(35).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 36 - SEP card service
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A36:D36").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(36, 1).Value = 44994
$ws.Cells.Item(36, 2).Value = 33
$ws.Cells.Item(36, 3).Value = 'Create SEP Card Service by using "ng g s ./@core/shared/services/sep-card"'
$ws.Cells.Item(36, 4).Value = "ng g s ./@core/shared/services/sep-card"

# ---------------------------------------------------------------------------
# Extend the table (ListObject) and its AutoFilter to cover the new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E36"))

# ---------------------------------------------------------------------------
# Update the saved view: scroll so row 28 is at the top and select D39,
# matching the state the workbook was left in.
# ---------------------------------------------------------------------------
$ws.Range("D39").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 28
